# question_5.xlsx edit script
# - Renames pivot-output sheet "Sheet2" -> "graph"
# - Removes the now-unused empty "Sheet1"
# - Corrects the GDPpct2 column: it should be last year's GDPpct (a one-year
#   lag of the GDPpct column), both in the source table (question_5) and in
#   the pivot-table output (graph / Sheet2), including the Grand Total row
# - Updates the pivot chart's series formulas to point at the renamed sheet
# - Refreshes window/selection state to match the edited workbook
# - Renames the chart shape to "graph_5"

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

$wsGraph = $wb.Worksheets.Item("Sheet2")
$wsData  = $wb.Worksheets.Item("question_5")

# ---------------------------------------------------------------------
# 1) Corrected GDPpct2 values (one-year lag of GDPpct), year 2007..2018
# ---------------------------------------------------------------------
$newGDPpct2 = @(6.77, 5.57, 2.69, -4.66, 2.4300000000000002, 1.76, -0.79, -0.05, 2.2599999999999998, 5.39, 2.54, 5.17)

# question_5 source table: years in rows 2..13, GDPpct2 in column J
for ($i = 0; $i -lt 12; $i++) {
    $row = 2 + $i
    $wsData.Cells.Item($row, 10).Value2 = $newGDPpct2[$i]
}

# graph (pivot output) sheet: years in rows 4..15, GDPpct2 in column C
for ($i = 0; $i -lt 12; $i++) {
    $row = 4 + $i
    $wsGraph.Cells.Item($row, 3).Value2 = $newGDPpct2[$i]
}

# Grand Total row (row 16) for the pivot output, column C
$total = 0
foreach ($v in $newGDPpct2) { $total = $total + $v }
$wsGraph.Cells.Item(16, 3).Value2 = $total

# ---------------------------------------------------------------------
# 2) Rename "Sheet2" -> "graph"
# ---------------------------------------------------------------------
$wsGraph.Name = "graph"

# ---------------------------------------------------------------------
# 3) Update the pivot chart series formulas to reference the new sheet name
# ---------------------------------------------------------------------
$co = $wsGraph.ChartObjects(1)
$co.Name = "graph_5"
$chart = $co.Chart
$chart.SeriesCollection(1).Formula = "=SERIES(graph!`$B`$3,graph!`$A`$4:`$A`$16,graph!`$B`$4:`$B`$16,1)"
$chart.SeriesCollection(2).Formula = "=SERIES(graph!`$C`$3,graph!`$A`$4:`$A`$16,graph!`$C`$4:`$C`$16,2)"
$chart.SeriesCollection(3).Formula = "=SERIES(graph!`$D`$3,graph!`$A`$4:`$A`$16,graph!`$D`$4:`$D`$16,3)"
$chart.SeriesCollection(4).Formula = "=SERIES(graph!`$E`$3,graph!`$A`$4:`$A`$16,graph!`$E`$4:`$E`$16,4)"

# ---------------------------------------------------------------------
# 4) Refresh pivot cache / pivot table so output reflects new data
# ---------------------------------------------------------------------
$wb.RefreshAll()
$pt = $wsGraph.PivotTables(1)
$pt.RefreshTable()

# ---------------------------------------------------------------------
# 5) Remove the empty, unused "Sheet1"
# ---------------------------------------------------------------------
$wsExtra = $wb.Worksheets.Item("Sheet1")
$wsExtra.Delete()

# ---------------------------------------------------------------------
# 6) Window / selection state
# ---------------------------------------------------------------------
$wsData.Select()
$winData = $excel.ActiveWindow
$winData.Zoom = 100
$wsData.Range("H12").Select()

$wsGraph.Select()
$winGraph = $excel.ActiveWindow
$winGraph.Zoom = 109
$wsGraph.Range("G33").Select()
